$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,23
$row2[0,0] = 1.885011672973633
$row2[0,1] = 2
$row2[0,2] = 6806.783671681492
$row2[0,3] = 0.220871380733333
$row2[0,4] = 0.2032341243744011
$row2[0,5] = 0.1811736602975099
$row2[0,6] = 0.1646373620654799
$row2[0,7] = 0.153424733464436
$row2[0,8] = 0.151810274894841
$row2[0,9] = 0.1500407159876893
$row2[0,10] = 0.1500407159876893
$row2[0,11] = 0.1498802964900214
$row2[0,12] = 0.1498802964900214
$row2[0,13] = 0.1498802964900214
$row2[0,14] = 0.1491726141855238
$row2[0,15] = 0.1487110274206565
$row2[0,16] = 0.1487110274206565
$row2[0,17] = 0.1487110274206565
$row2[0,18] = 0.1487110274206565
$row2[0,19] = 0.1487110274206565
$row2[0,20] = 0.1487110274206565
$row2[0,21] = 0.1487110274206565
$row2[0,22] = 0.1486858415532454
$ws.Range("C2:Y2").Value = $row2

$ws.Range("C3").Value = 1.623996734619141
$row3 = New-Object 'object[,]' 1,21
$row3[0,0] = 6884.128429302994
$row3[0,1] = 0.220871380733333
$row3[0,2] = 0.2032341243744011
$row3[0,3] = 0.1872033993012155
$row3[0,4] = 0.1701453130149286
$row3[0,5] = 0.1601306333182273
$row3[0,6] = 0.1519673006992826
$row3[0,7] = 0.1505809403848549
$row3[0,8] = 0.1505809403848549
$row3[0,9] = 0.1505809403848549
$row3[0,10] = 0.1505809403848549
$row3[0,11] = 0.1505809403848549
$row3[0,12] = 0.1504224749945868
$row3[0,13] = 0.1503328034145848
$row3[0,14] = 0.1503328034145848
$row3[0,15] = 0.1503328034145848
$row3[0,16] = 0.1503328034145848
$row3[0,17] = 0.1503328034145848
$row3[0,18] = 0.1503328034145848
$row3[0,19] = 0.1502721948985342
$row3[0,20] = 0.1501935366335866
$ws.Range("E3:Y3").Value = $row3

$ws.Range("C4").Value = 1.711998224258423
$row4 = New-Object 'object[,]' 1,21
$row4[0,0] = 6801.461594161514
$row4[0,1] = 0.220871380733333
$row4[0,2] = 0.2032341243744011
$row4[0,3] = 0.181428182682164
$row4[0,4] = 0.1656014597423033
$row4[0,5] = 0.1573990214038638
$row4[0,6] = 0.1506633937187187
$row4[0,7] = 0.1489662791034382
$row4[0,8] = 0.1489662791034382
$row4[0,9] = 0.1489662791034382
$row4[0,10] = 0.1489662791034382
$row4[0,11] = 0.1489662791034382
$row4[0,12] = 0.1489662791034382
$row4[0,13] = 0.1488605354935528
$row4[0,14] = 0.1488605354935528
$row4[0,15] = 0.1488605354935528
$row4[0,16] = 0.1488605354935528
$row4[0,17] = 0.1488605354935528
$row4[0,18] = 0.1487250669840534
$row4[0,19] = 0.1487250669840534
$row4[0,20] = 0.1485820973520763
$ws.Range("E4:Y4").Value = $row4

$ws.Range("C5").Value = 1.66700267791748
$row5 = New-Object 'object[,]' 1,21
$row5[0,0] = 6712.782018365953
$row5[0,1] = 0.220871380733333
$row5[0,2] = 0.2032341243744011
$row5[0,3] = 0.1906308968166822
$row5[0,4] = 0.169904710840951
$row5[0,5] = 0.1537924725610029
$row5[0,6] = 0.1513749513919128
$row5[0,7] = 0.1513749513919128
$row5[0,8] = 0.148322504744229
$row5[0,9] = 0.147657906098213
$row5[0,10] = 0.147657906098213
$row5[0,11] = 0.1476018783594891
$row5[0,12] = 0.1476018783594891
$row5[0,13] = 0.1476018783594891
$row5[0,14] = 0.1475136820048581
$row5[0,15] = 0.1475136820048581
$row5[0,16] = 0.1468534506504084
$row5[0,17] = 0.1468534506504084
$row5[0,18] = 0.1468534506504084
$row5[0,19] = 0.1468534506504084
$row5[0,20] = 0.1468534506504084
$ws.Range("E5:Y5").Value = $row5

$ws.Range("C6").Value = 1.814995765686035
$row6 = New-Object 'object[,]' 1,21
$row6[0,0] = 6925.030223933425
$row6[0,1] = 0.220871380733333
$row6[0,2] = 0.2032341243744011
$row6[0,3] = 0.1744302171245432
$row6[0,4] = 0.1716635861134211
$row6[0,5] = 0.1628475448196577
$row6[0,6] = 0.1628475448196577
$row6[0,7] = 0.1594992282521263
$row6[0,8] = 0.1594992282521263
$row6[0,9] = 0.1585767820002412
$row6[0,10] = 0.1579806163026953
$row6[0,11] = 0.1521911351186682
$row6[0,12] = 0.1513932547330028
$row6[0,13] = 0.1513932547330028
$row6[0,14] = 0.1513932547330028
$row6[0,15] = 0.1513932547330028
$row6[0,16] = 0.1513932547330028
$row6[0,17] = 0.1511459024096665
$row6[0,18] = 0.1511459024096665
$row6[0,19] = 0.1509979260220096
$row6[0,20] = 0.1509908425718016
$ws.Range("E6:Y6").Value = $row6

$ws.Range("C7").Value = 1.653035879135132
$row7 = New-Object 'object[,]' 1,21
$row7[0,0] = 6757.639205610843
$row7[0,1] = 0.220871380733333
$row7[0,2] = 0.2032341243744011
$row7[0,3] = 0.1901283472164353
$row7[0,4] = 0.1748017811196545
$row7[0,5] = 0.1586952293583923
$row7[0,6] = 0.1543582244312179
$row7[0,7] = 0.1517177536499814
$row7[0,8] = 0.150174287939802
$row7[0,9] = 0.1482042938436714
$row7[0,10] = 0.1482042938436714
$row7[0,11] = 0.1482042938436714
$row7[0,12] = 0.1482042938436714
$row7[0,13] = 0.1481233185688406
$row7[0,14] = 0.1481233185688406
$row7[0,15] = 0.1481233185688406
$row7[0,16] = 0.1481233185688406
$row7[0,17] = 0.1481233185688406
$row7[0,18] = 0.1479637093993097
$row7[0,19] = 0.1479637093993097
$row7[0,20] = 0.1477278597584959
$ws.Range("E7:Y7").Value = $row7

$ws.Range("C8").Value = 1.774947643280029
$row8 = New-Object 'object[,]' 1,21
$row8[0,0] = 6692.330289027674
$row8[0,1] = 0.220871380733333
$row8[0,2] = 0.2032341243744011
$row8[0,3] = 0.186494837268027
$row8[0,4] = 0.1639015292423214
$row8[0,5] = 0.1590414918004688
$row8[0,6] = 0.1501816870455802
$row8[0,7] = 0.1484580760284703
$row8[0,8] = 0.1484580760284703
$row8[0,9] = 0.1480200015087329
$row8[0,10] = 0.1474528095956709
$row8[0,11] = 0.1474528095956709
$row8[0,12] = 0.1474528095956709
$row8[0,13] = 0.1474528095956709
$row8[0,14] = 0.1471905636772994
$row8[0,15] = 0.1471905636772994
$row8[0,16] = 0.1471905636772994
$row8[0,17] = 0.1468482393236899
$row8[0,18] = 0.1465003364333648
$row8[0,19] = 0.1465003364333648
$row8[0,20] = 0.1464547814625277
$ws.Range("E8:Y8").Value = $row8

$ws.Range("C9").Value = 1.58800220489502
$row9 = New-Object 'object[,]' 1,21
$row9[0,0] = 6904.865481738917
$row9[0,1] = 0.220871380733333
$row9[0,2] = 0.2032341243744011
$row9[0,3] = 0.1695993794030745
$row9[0,4] = 0.1541227023999988
$row9[0,5] = 0.1524378853406889
$row9[0,6] = 0.1519122789697667
$row9[0,7] = 0.1505977676752225
$row9[0,8] = 0.1505977676752225
$row9[0,9] = 0.1505977676752225
$row9[0,10] = 0.1505977676752225
$row9[0,11] = 0.1505977676752225
$row9[0,12] = 0.1505977676752225
$row9[0,13] = 0.1505977676752225
$row9[0,14] = 0.1505977676752225
$row9[0,15] = 0.1505977676752225
$row9[0,16] = 0.1505977676752225
$row9[0,17] = 0.1505977676752225
$row9[0,18] = 0.1505977676752225
$row9[0,19] = 0.1505977676752225
$row9[0,20] = 0.1505977676752225
$ws.Range("E9:Y9").Value = $row9

$ws.Range("C10").Value = 1.646013498306274
$row10 = New-Object 'object[,]' 1,21
$row10[0,0] = 6682.48225932294
$row10[0,1] = 0.220871380733333
$row10[0,2] = 0.2032341243744011
$row10[0,3] = 0.1868489508535622
$row10[0,4] = 0.1726400685735788
$row10[0,5] = 0.1532344469202108
$row10[0,6] = 0.1488140718027697
$row10[0,7] = 0.1488140718027697
$row10[0,8] = 0.1488140718027697
$row10[0,9] = 0.1488140718027697
$row10[0,10] = 0.1487365734560907
$row10[0,11] = 0.1486382911904435
$row10[0,12] = 0.1474241066977662
$row10[0,13] = 0.1474241066977662
$row10[0,14] = 0.1469875173914951
$row10[0,15] = 0.1467019731976454
$row10[0,16] = 0.1464912416848005
$row10[0,17] = 0.146402187889474
$row10[0,18] = 0.1462819105192005
$row10[0,19] = 0.1462819105192005
$row10[0,20] = 0.1462628120725719
$ws.Range("E10:Y10").Value = $row10

$ws.Range("C11").Value = 1.65595817565918
$row11 = New-Object 'object[,]' 1,21
$row11[0,0] = 6787.794216585053
$row11[0,1] = 0.220871380733333
$row11[0,2] = 0.2032341243744011
$row11[0,3] = 0.1738589344658298
$row11[0,4] = 0.1679048363105065
$row11[0,5] = 0.1595504085035202
$row11[0,6] = 0.1535498982761371
$row11[0,7] = 0.1535498982761371
$row11[0,8] = 0.1526275233419406
$row11[0,9] = 0.15108600040064
$row11[0,10] = 0.1496161460055117
$row11[0,11] = 0.1494513086152274
$row11[0,12] = 0.1489300881227507
$row11[0,13] = 0.1489300881227507
$row11[0,14] = 0.1489300881227507
$row11[0,15] = 0.1489300881227507
$row11[0,16] = 0.1489300881227507
$row11[0,17] = 0.1484224038657344
$row11[0,18] = 0.1484224038657344
$row11[0,19] = 0.1484224038657344
$row11[0,20] = 0.1483156767365507
$ws.Range("E11:Y11").Value = $row11

Write-Host "done"